# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets
# to match newly scraped totals (gh-pages data output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 520
$ws1.Range("F8").Value  = 493
$ws1.Range("F9").Value  = 6506
$ws1.Range("F10").Value = 174
$ws1.Range("F11").Value = 136
$ws1.Range("F12").Value = 1022
$ws1.Range("F13").Value = 335
$ws1.Range("F14").Value = 105
$ws1.Range("F15").Value = 179
$ws1.Range("F16").Value = 482

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 520
$ws4.Range("F8").Value  = 493
$ws4.Range("F9").Value  = 6506
$ws4.Range("F10").Value = 174
$ws4.Range("F11").Value = 136
$ws4.Range("F12").Value = 1022
$ws4.Range("F13").Value = 336
$ws4.Range("F14").Value = 105
$ws4.Range("F15").Value = 179
$ws4.Range("F16").Value = 482
